$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

# Update the product identifier text on both sheets (dash added after "342")
$newName = "342-MS-EPP-DB-SAR-REC-NON-RNI-CTPD-DL-MD-TR-1-ONTIME"
$wsInput.Range("B1").Value = $newName
$wsOutput.Range("B1").Value = $newName

# Update selection on the input sheet (was A6:B6, now B1)
$wsInput.Range("B1").Select() | Out-Null

# Update selection on the output sheet and make it the active tab
$wsOutput.Select() | Out-Null
$wsOutput.Range("B1").Select() | Out-Null
